$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Malik Monk", "PG,SG,SF", "Sacramento Kings"),
    @("Andrew Wiggins", "SF,PF", "Miami Heat"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Keon Ellis", "SG,SF", "Sacramento Kings"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Cade Cunningham", "PG,SG", "Detroit Pistons"),
    @("Carlton Carrington", "PG,SG", "Washington Wizards"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Ausar Thompson", "SF,PF", "Detroit Pistons"),
    @("Precious Achiuwa", "PF,C", "New York Knicks"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Damian Lillard", "PG", "Milwaukee Bucks"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Anthony Davis", "PF,C", "Dallas Mavericks")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row = $row + 1
}
